# Scheduled runner update: refresh market-price / profit figures in
# Kujata_Profits.xlsx (currentAveragePrice*, LevePrice*, LeveProfit* columns
# H:N) for a handful of Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/WVR
# sheets, as pulled from the latest Universalis snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 136
$ws.Range("I9").Value = 120
$ws.Range("K9").Value = 120
$ws.Range("M9").Value = 49

$ws.Range("H12").Value = 406.125
$ws.Range("I12").Value = 433
$ws.Range("J12").Value = 390
$ws.Range("K12").Value = 433
$ws.Range("L12").Value = 390
$ws.Range("M12").Value = -263
$ws.Range("N12").Value = -730

$ws.Range("H15").Value = 2857.1294
$ws.Range("I15").Value = 2857.1294
$ws.Range("K15").Value = 8571.388199999999
$ws.Range("M15").Value = -8402.388199999999

$ws.Range("H21").Value = 29919
$ws.Range("J21").Value = 29919
$ws.Range("L21").Value = 29919
$ws.Range("N21").Value = -30855

$ws.Range("H23").Value = 29919
$ws.Range("J23").Value = 29919
$ws.Range("L23").Value = 29919
$ws.Range("N23").Value = -30387

$ws.Range("H29").Value = 1821.0526
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 2253.3333
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 6759.999899999999
$ws.Range("M29").Value = -319
$ws.Range("N29").Value = -7321.999899999999

$ws.Range("H38").Value = 1628.975
$ws.Range("J38").Value = 1995.1875
$ws.Range("L38").Value = 5985.5625
$ws.Range("N38").Value = -6729.5625

$ws.Range("H43").Value = 7952793.5
$ws.Range("I43").Value = 50500.5
$ws.Range("J43").Value = 11113711
$ws.Range("K43").Value = 50500.5
$ws.Range("L43").Value = 11113711
$ws.Range("M43").Value = -50431.5
$ws.Range("N43").Value = -11113849

$ws.Range("H58").Value = 1181.6111
$ws.Range("I58").Value = 455.16666
$ws.Range("J58").Value = 2634.5
$ws.Range("K58").Value = 1365.49998
$ws.Range("L58").Value = 7903.5
$ws.Range("M58").Value = -1215.49998
$ws.Range("N58").Value = -8203.5

$ws.Range("H132").Value = 8842.709999999999
$ws.Range("I132").Value = 7288.923
$ws.Range("K132").Value = 21866.769
$ws.Range("M132").Value = -19336.769

$ws.Range("H134").Value = 37158.8
$ws.Range("J134").Value = 37158.8
$ws.Range("L134").Value = 37158.8
$ws.Range("N134").Value = -47298.8

$ws.Range("H141").Value = 1148.4615
$ws.Range("I141").Value = 994.1667
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 2982.5001
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 2197.4999
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1671.3334
$ws.Range("I45").Value = 1500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1123

$ws.Range("H132").Value = 3421.6667
$ws.Range("I132").Value = 3172.3333
$ws.Range("K132").Value = 9516.999899999999
$ws.Range("M132").Value = -6986.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24472.428
$ws.Range("I82").Value = 8153.5
$ws.Range("J82").Value = 31000
$ws.Range("K82").Value = 8153.5
$ws.Range("L82").Value = 31000
$ws.Range("M82").Value = -7770.5
$ws.Range("N82").Value = -31766

$ws.Range("H85").Value = 24472.428
$ws.Range("I85").Value = 8153.5
$ws.Range("J85").Value = 31000
$ws.Range("K85").Value = 8153.5
$ws.Range("L85").Value = 31000
$ws.Range("M85").Value = -6827.5
$ws.Range("N85").Value = -33652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1728.6072
$ws.Range("J31").Value = 1443.1765
$ws.Range("L31").Value = 1443.1765
$ws.Range("N31").Value = -2033.1765

$ws.Range("H34").Value = 1728.6072
$ws.Range("J34").Value = 1443.1765
$ws.Range("L34").Value = 1443.1765
$ws.Range("N34").Value = -1847.1765

$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496

$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716

$ws.Range("H134").Value = 15153328
$ws.Range("I134").Value = 1833.9
$ws.Range("K134").Value = 5501.700000000001
$ws.Range("M134").Value = -2966.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 890.9091
$ws.Range("I17").Value = 350
$ws.Range("J17").Value = 1200
$ws.Range("K17").Value = 1050
$ws.Range("L17").Value = 3600
$ws.Range("M17").Value = -881
$ws.Range("N17").Value = -3938

$ws.Range("H34").Value = 2313.7693
$ws.Range("J34").Value = 2689
$ws.Range("L34").Value = 8067
$ws.Range("N34").Value = -8235

$ws.Range("H55").Value = 1835.4286
$ws.Range("J55").Value = 2065.6667
$ws.Range("L55").Value = 6197.000100000001
$ws.Range("N55").Value = -6551.000100000001

$ws.Range("J131").Value = 2304.6956
$ws.Range("L131").Value = 6914.0868
$ws.Range("N131").Value = -16994.0868

$ws.Range("H133").Value = 4334.5264
$ws.Range("I133").Value = 3093.3333
$ws.Range("J133").Value = 4907.385
$ws.Range("K133").Value = 9279.999899999999
$ws.Range("L133").Value = 14722.155
$ws.Range("M133").Value = -4219.999899999999
$ws.Range("N133").Value = -24842.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5200
$ws.Range("I80").Value = 4625
$ws.Range("J80").Value = 5966.6665
$ws.Range("K80").Value = 4625
$ws.Range("L80").Value = 5966.6665
$ws.Range("M80").Value = -3627
$ws.Range("N80").Value = -7962.6665

$ws.Range("H83").Value = 5200
$ws.Range("I83").Value = 4625
$ws.Range("J83").Value = 5966.6665
$ws.Range("K83").Value = 23125
$ws.Range("L83").Value = 29833.3325
$ws.Range("M83").Value = -18133
$ws.Range("N83").Value = -39817.3325

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 250.16667
$ws.Range("I81").Value = 250.16667
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 500.33334
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 560.66666
$ws.Range("N81").Value = $null

$ws.Range("H84").Value = 250.16667
$ws.Range("I84").Value = 250.16667
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 2501.6667
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 2802.3333
$ws.Range("N84").Value = $null
